$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.249.42'
$ws.Range("E2").Value = '  +1.05%  '

$ws.Range("D3").Value = '1.852.40'
$ws.Range("E3").Value = '  +1.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.93'
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("E6").Value = '  -0.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4600'
$ws.Range("E7").Value = '  +0.85%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3707'
$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07291'
$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8866'
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.09'
$ws.Range("E11").Value = '  +1.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07827'
$ws.Range("E12").Value = '  -1.69%  '

$ws.Range("D13").Value = '1.807.34'
$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.388'
$ws.Range("E14").Value = '  +1.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.527'
$ws.Range("E15").Value = '  +0.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.40'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.76'
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("D21").Value = '27.277.61'
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.110'
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").Value = '2.054.58'
$ws.Range("E24").Value = '  +3.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.916'
$ws.Range("E25").Value = '  +4.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.95'
$ws.Range("E26").Value = '  -0.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.43'
$ws.Range("E27").Value = '  +0.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.057'
$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.81'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.061'
$ws.Range("E30").Value = '  -1.58%  '

$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7729'
$ws.Range("E32").Value = '  +6.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.074'
$ws.Range("E33").Value = '  +3.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.168'
$ws.Range("E34").Value = '  +3.59%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.495'
$ws.Range("E35").Value = '  +2.05%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.745'
$ws.Range("E36").Value = '  +12.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.084'
$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05260'
$ws.Range("E39").Value = '  +0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.952'
$ws.Range("E40").Value = '  +0.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.057'
$ws.Range("E41").Value = '  -0.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5120'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1634'
$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.391'
$ws.Range("E44").Value = '  +2.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4795'
$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.34'
$ws.Range("E46").Value = '  +1.27%  '

$ws.Range("E47").Value = '  -0.57%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.14'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.642'
$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.64'
$ws.Range("E51").Value = '  +0.68%  '
